# Weekly update: insert a new price-report row for "Cilantro" at
# Vega Monumental Concepción, shifting the existing rows 129-218 down by
# one (to 130-219) and populating the newly opened row 129 with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 129 - Excel shifts rows 129:218 down to
# 130:219 and carries the existing formatting (e.g. the date style on
# column D) down with them.
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with this week's record.
$ws.Cells.Item(129, 1).Value2  = 11
$ws.Cells.Item(129, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(129, 3).Value2  = "Bíobío"
$ws.Cells.Item(129, 4).Value2  = 44777
$ws.Cells.Item(129, 5).Value2  = 8
$ws.Cells.Item(129, 6).Value2  = 100112040
$ws.Cells.Item(129, 7).Value2  = "Cilantro"
$ws.Cells.Item(129, 8).Value2  = "Sin especificar"
$ws.Cells.Item(129, 9).Value2  = "Primera"
$ws.Cells.Item(129, 10).Value2 = 120
$ws.Cells.Item(129, 11).Value2 = 8000
$ws.Cells.Item(129, 12).Value2 = 8500
$ws.Cells.Item(129, 13).Value2 = 8292
$ws.Cells.Item(129, 14).Value2 = "$/caja 36 atados"
$ws.Cells.Item(129, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(129, 16).Value2 = 230
$ws.Cells.Item(129, 17).Value2 = 36
$ws.Cells.Item(129, 18).Value2 = "Hortaliza"
